$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "Friendly Matches / Tombense MG vs Desportiva" row.
# This shifts all subsequent rows up by one (old row 3 -> new row 2, etc.).
$ws.Rows(2).Delete()

# Refresh the odds data for the remaining six fixtures (now rows 2-7).
$data = New-Object 'object[,]' 6,41

# Row 2: MC Alger vs ES Ben Aknoun (Algerian Ligue 1)
$data[0,0] = 'Algerian Ligue 1'
$data[0,1] = '2025-12-23'
$data[0,2] = '15:30:00'
$data[0,3] = 'MC Alger'
$data[0,4] = 'ES Ben Aknoun'
$data[0,5] = 1.37
$data[0,6] = 1.41
$data[0,7] = 13
$data[0,8] = 23
$data[0,9] = 4.5
$data[0,10] = 5.4
$data[0,11] = 1.5
$data[0,12] = 1.09
$data[0,13] = 3.05
$data[0,14] = 1.43
$data[0,15] = 1.71
$data[0,16] = 2.24
$data[0,17] = 1.25
$data[0,18] = 4.4
$data[0,19] = 2.74
$data[0,20] = 1.54
$data[0,21] = 1.07
$data[0,22] = 3.3
$data[0,23] = 1000
$data[0,24] = 30
$data[0,25] = 1000
$data[0,26] = 1000
$data[0,27] = 6.6
$data[0,28] = 970
$data[0,29] = 1000
$data[0,30] = 1000
$data[0,31] = 7.6
$data[0,32] = 16
$data[0,33] = 1000
$data[0,34] = 1000
$data[0,35] = 900
$data[0,36] = 29
$data[0,37] = 1000
$data[0,38] = 1000
$data[0,39] = 12.5
$data[0,40] = 1000

# Row 3: Academia de Balompie Boliviano vs San Juan FC (Bolivian Liga de Futbol Profesional)
$data[1,0] = 'Bolivian Liga de Futbol Profesional'
$data[1,1] = '2025-12-23'
$data[1,2] = '16:00:00'
$data[1,3] = 'Academia de Balompie Boliviano'
$data[1,4] = 'San Juan FC'
$data[1,5] = 1.33
$data[1,6] = 1.35
$data[1,7] = 8
$data[1,8] = 11
$data[1,9] = 6.2
$data[1,10] = 9.199999999999999
$data[1,11] = 1.25
$data[1,12] = 1.03
$data[1,13] = 6.4
$data[1,14] = 1.14
$data[1,15] = 3
$data[1,16] = 1.45
$data[1,17] = 1.79
$data[1,18] = 2.12
$data[1,19] = 1.84
$data[1,20] = 1.84
$data[1,21] = 1.11
$data[1,22] = 3.25
$data[1,23] = 1000
$data[1,24] = 1000
$data[1,25] = 1000
$data[1,26] = 1000
$data[1,27] = 1000
$data[1,28] = 1000
$data[1,29] = 1000
$data[1,30] = 1000
$data[1,31] = 1000
$data[1,32] = 1000
$data[1,33] = 1000
$data[1,34] = 1000
$data[1,35] = 1000
$data[1,36] = 1000
$data[1,37] = 1000
$data[1,38] = 1000
$data[1,39] = 1000
$data[1,40] = 1000

# Row 4: Serra Branca EC vs Maguary (Friendly Matches)
$data[2,0] = 'Friendly Matches'
$data[2,1] = '2025-12-23'
$data[2,2] = '16:00:00'
$data[2,3] = 'Serra Branca EC'
$data[2,4] = 'Maguary'
$data[2,5] = 1.64
$data[2,6] = 2.2
$data[2,7] = 2.8
$data[2,8] = 7.2
$data[2,9] = 2.92
$data[2,10] = 7.2
$data[2,11] = 1.36
$data[2,12] = 1.07
$data[2,13] = 1.1
$data[2,14] = 1.29
$data[2,15] = 1.6
$data[2,16] = 1.5
$data[2,17] = 1.41
$data[2,18] = 2.4
$data[2,19] = 1.05
$data[2,20] = 1.04
$data[2,21] = 1.25
$data[2,22] = 1.83
$data[2,23] = 990
$data[2,24] = 1000
$data[2,25] = 1000
$data[2,26] = 1000
$data[2,27] = 11
$data[2,28] = 19
$data[2,29] = 1000
$data[2,30] = 1000
$data[2,31] = 1000
$data[2,32] = 40
$data[2,33] = 60
$data[2,34] = 1000
$data[2,35] = 1000
$data[2,36] = 1000
$data[2,37] = 1000
$data[2,38] = 1000
$data[2,39] = 1000
$data[2,40] = 1000

# Row 5: Guimaraes vs Sporting Lisbon (Portuguese Primeira Liga)
$data[3,0] = 'Portuguese Primeira Liga'
$data[3,1] = '2025-12-23'
$data[3,2] = '17:45:00'
$data[3,3] = 'Guimaraes'
$data[3,4] = 'Sporting Lisbon'
$data[3,5] = 8
$data[3,6] = 8.199999999999999
$data[3,7] = 1.47
$data[3,8] = 1.49
$data[3,9] = 4.8
$data[3,10] = 5
$data[3,11] = 1.39
$data[3,12] = 1.06
$data[3,13] = 4.1
$data[3,14] = 1.3
$data[3,15] = 2.08
$data[3,16] = 1.89
$data[3,17] = 1.4
$data[3,18] = 3.4
$data[3,19] = 2.12
$data[3,20] = 1.86
$data[3,21] = 3
$data[3,22] = 1.13
$data[3,23] = 17.5
$data[3,24] = 7.8
$data[3,25] = 8.199999999999999
$data[3,26] = 12.5
$data[3,27] = 24
$data[3,28] = 10.5
$data[3,29] = 9.6
$data[3,30] = 15.5
$data[3,31] = 70
$data[3,32] = 32
$data[3,33] = 27
$data[3,34] = 40
$data[3,35] = 280
$data[3,36] = 140
$data[3,37] = 120
$data[3,38] = 160
$data[3,39] = 160
$data[3,40] = 7.6

# Row 6: Necaxa vs Atletico San Luis (Friendly Matches)
$data[4,0] = 'Friendly Matches'
$data[4,1] = '2025-12-23'
$data[4,2] = '18:00:00'
$data[4,3] = 'Necaxa'
$data[4,4] = 'Atletico San Luis'
$data[4,5] = 1.86
$data[4,6] = 2.32
$data[4,7] = 3.25
$data[4,8] = 4.9
$data[4,9] = 3
$data[4,10] = 4.2
$data[4,11] = 1.34
$data[4,12] = 1.07
$data[4,13] = 3.05
$data[4,14] = 1.3
$data[4,15] = 1.94
$data[4,16] = 1.47
$data[4,17] = 1.36
$data[4,18] = 2.36
$data[4,19] = 1.04
$data[4,20] = 1.04
$data[4,21] = 1.25
$data[4,22] = 1.76
$data[4,23] = 1000
$data[4,24] = 1000
$data[4,25] = 1000
$data[4,26] = 1000
$data[4,27] = 46
$data[4,28] = 19
$data[4,29] = 1000
$data[4,30] = 1000
$data[4,31] = 1000
$data[4,32] = 40
$data[4,33] = 990
$data[4,34] = 1000
$data[4,35] = 1000
$data[4,36] = 1000
$data[4,37] = 1000
$data[4,38] = 1000
$data[4,39] = 85
$data[4,40] = 1000

# Row 7: Real Espana vs CD Motagua (Honduras Liga Nacional)
$data[5,0] = 'Honduras Liga Nacional'
$data[5,1] = '2025-12-23'
$data[5,2] = '22:00:00'
$data[5,3] = 'Real Espana'
$data[5,4] = 'CD Motagua'
$data[5,5] = 1.76
$data[5,6] = 1.84
$data[5,7] = 4.9
$data[5,8] = 6.4
$data[5,9] = 3.45
$data[5,10] = 3.95
$data[5,11] = 1.43
$data[5,12] = 1.07
$data[5,13] = 3.5
$data[5,14] = 1.33
$data[5,15] = 1.87
$data[5,16] = 1.95
$data[5,17] = 1.32
$data[5,18] = 3.4
$data[5,19] = 1.87
$data[5,20] = 1.86
$data[5,21] = 1.2
$data[5,22] = 2.18
$data[5,23] = 1000
$data[5,24] = 18.5
$data[5,25] = 42
$data[5,26] = 160
$data[5,27] = 1000
$data[5,28] = 1000
$data[5,29] = 23
$data[5,30] = 85
$data[5,31] = 10.5
$data[5,32] = 10
$data[5,33] = 24
$data[5,34] = 1000
$data[5,35] = 19
$data[5,36] = 20
$data[5,37] = 40
$data[5,38] = 140
$data[5,39] = 1000
$data[5,40] = 110

$ws.Range("A2:AO7").Value = $data

# Column B ("Date") holds plain text like "2025-12-23"; a bare assignment would
# otherwise be auto-recognized as a date and stored as a serial number. Force
# it to remain text, then clear the temporary number format back to the default
# "Normal" style so no extra formatting is introduced on the cells.
$dateRange = $ws.Range("B2:B7")
$dateRange.NumberFormat = "@"
$dateRange.Value = "2025-12-23"
$dateRange.Style = "Normal"

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
